$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New query text (CDS Study filter fixes): all three prior queries (Participants,
# Stats, Files) are rewritten, and the Samples query gets an updated WHERE clause.
# Cypher backticks are doubled (``) to emit a literal backtick in PowerShell
# here-strings; double/single quotes pass through unescaped.
# ---------------------------------------------------------------------------

$queryParticipants = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$querySamples = @"
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as ``Sample ID``,
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(samp.sample_tumor_status,'') as ``Tumor``,
coalesce(samp.sample_type,'') as ``Analyte Type``
  ORDER By samp.sample_id 
  LIMIT 100
"@

$queryFiles = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name
Limit 100
"@

$queryStats = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

# ---------------------------------------------------------------------------
# Row 2 now describes the ParticipantsTab (previously row 2 held the old
# Participants-style query under the "ParticipantsTab" label already, so the
# TabName column stays "ParticipantsTab"; only the query text changes).
# Row 3 keeps SamplesTab, row 4 keeps FilesTab.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $queryParticipants
$ws.Range("C2").Value = $queryStats

$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $querySamples
$ws.Range("C3").Value = $queryStats

$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $queryFiles
$ws.Range("C4").Value = $queryStats

# ---------------------------------------------------------------------------
# Column A widens from bestFit ~11.57 to a fixed custom width of 19 characters.
# (18.15 is the COM "points" input that this runtime's px-rounding resolves to
# an OOXML character width of exactly 19.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.15

# ---------------------------------------------------------------------------
# Rows 2 and 4 grow to Excel's maximum row height because their query text got
# much longer; row 3's height is unaffected by the edit, so it is pinned back
# to its original value (writing the long StatQuery text into column C makes
# this runtime provisionally grow row 3's height, which we then undo).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 218.25

# Selection moves from A2 to B2.
$ws.Range("B2").Select()
